$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Range("H8").Value = 749.6667
$ws.Range("I8").Value = 99.59999999999999
$ws.Range("J8").Value = 4000
$ws.Range("K8").Value = 298.8
$ws.Range("L8").Value = 12000
$ws.Range("M8").Value = -159.8
$ws.Range("N8").Value = -12278
# Row 40
$ws.Range("H40").Value = 3365
$ws.Range("I40").Value = 2750
$ws.Range("J40").Value = 3980
$ws.Range("K40").Value = 2750
$ws.Range("L40").Value = 3980
$ws.Range("M40").Value = -2575
$ws.Range("N40").Value = -4330
# Row 127
$ws.Range("H127").Value = 1969.1875
$ws.Range("J127").Value = 2571.5715
$ws.Range("L127").Value = 7714.7145
$ws.Range("N127").Value = -17634.7145
# Row 131
$ws.Range("H131").Value = 696.4286
$ws.Range("I131").Value = 696.4286
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 2089.2858
$ws.Range("L131").Value = 0
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = 2950.7142
# Row 138
$ws.Range("H138").Value = 2057.9
$ws.Range("I138").Value = 1541.3636
$ws.Range("J138").Value = 2463.75
$ws.Range("K138").Value = 4624.0908
$ws.Range("L138").Value = 7391.25
$ws.Range("M138").Value = 515.9092000000001
$ws.Range("N138").Value = -17671.25
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 929.3182
$ws.Range("I2").Value = 913.7143
$ws.Range("J2").Value = 956.625
$ws.Range("K2").Value = 913.7143
$ws.Range("L2").Value = 956.625
$ws.Range("M2").Value = -800.7143
$ws.Range("N2").Value = -1182.625
# Row 45
$ws.Range("H45").Value = 1290.2667
$ws.Range("I45").Value = 1671.75
$ws.Range("J45").Value = 854.2857
$ws.Range("K45").Value = 1671.75
$ws.Range("L45").Value = 854.2857
$ws.Range("M45").Value = -1294.75
$ws.Range("N45").Value = -1608.2857
# Row 82
$ws.Range("H82").Value = 37400
$ws.Range("I82").Value = 30000
$ws.Range("J82").Value = 39866.668
$ws.Range("K82").Value = 30000
$ws.Range("L82").Value = 39866.668
$ws.Range("M82").Value = -29639
$ws.Range("N82").Value = -40588.668
# Row 85
$ws.Range("H85").Value = 37400
$ws.Range("I85").Value = 30000
$ws.Range("J85").Value = 39866.668
$ws.Range("K85").Value = 30000
$ws.Range("L85").Value = 39866.668
$ws.Range("M85").Value = -28752
$ws.Range("N85").Value = -42362.668
# Row 110
$ws.Range("H110").Value = 828.1177
$ws.Range("I110").Value = 763.9
$ws.Range("J110").Value = 919.8570999999999
$ws.Range("K110").Value = 763.9
$ws.Range("L110").Value = 919.8570999999999
$ws.Range("M110").Value = 1281.1
$ws.Range("N110").Value = -5009.8571
# Row 116
$ws.Range("H116").Value = 929.3182
$ws.Range("I116").Value = 913.7143
$ws.Range("J116").Value = 956.625
$ws.Range("K116").Value = 913.7143
$ws.Range("L116").Value = 956.625
$ws.Range("M116").Value = 1380.2857
$ws.Range("N116").Value = -5544.625
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 929.3182
$ws.Range("I3").Value = 913.7143
$ws.Range("J3").Value = 956.625
$ws.Range("K3").Value = 913.7143
$ws.Range("L3").Value = 956.625
$ws.Range("M3").Value = -799.7143
$ws.Range("N3").Value = -1184.625
# Row 58
$ws.Range("H58").Value = 30340.666
$ws.Range("J58").Value = 30340.666
$ws.Range("L58").Value = 30340.666
$ws.Range("N58").Value = -30928.666
# Row 59
$ws.Range("H59").Value = 61833.75
$ws.Range("J59").Value = 61833.75
$ws.Range("L59").Value = 61833.75
$ws.Range("N59").Value = -63527.75
# Row 105
$ws.Range("H105").Value = 2102.75
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()
# Row 107
$ws.Range("H107").Value = 90911130
$ws.Range("I107").Value = 142859650
$ws.Range("J107").Value = 1200
$ws.Range("K107").Value = 142859650
$ws.Range("L107").Value = 1200
$ws.Range("M107").Value = -142857730
$ws.Range("N107").Value = -5040
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 6264363
$ws.Range("I58").Value = 8991940
$ws.Range("J58").Value = 29902
$ws.Range("K58").Value = 8991940
$ws.Range("L58").Value = 29902
$ws.Range("M58").Value = -8991737
$ws.Range("N58").Value = -30308
# Row 122
$ws.Range("H122").Value = 4330569
$ws.Range("I122").Value = 7143891
$ws.Range("J122").Value = 2381.8462
$ws.Range("K122").Value = 21431673
$ws.Range("L122").Value = 7145.5386
$ws.Range("M122").Value = -21429223
$ws.Range("N122").Value = -12045.5386
# Row 132
$ws.Range("H132").Value = 8337477
$ws.Range("I132").Value = 15873893
$ws.Range("J132").Value = 7754.316
$ws.Range("K132").Value = 47621679
$ws.Range("L132").Value = 23262.948
$ws.Range("M132").Value = -47619149
$ws.Range("N132").Value = -28322.948
# Row 134
$ws.Range("H134").Value = 9470915
$ws.Range("I134").Value = 8621766
$ws.Range("J134").Value = 15627250
$ws.Range("K134").Value = 25865298
$ws.Range("L134").Value = 46881750
$ws.Range("M134").Value = -25862763
$ws.Range("N134").Value = -46886820
# Row 136
$ws.Range("H136").Value = 6264363
$ws.Range("I136").Value = 8991940
$ws.Range("J136").Value = 29902
$ws.Range("K136").Value = 26975820
$ws.Range("L136").Value = 89706
$ws.Range("M136").Value = -26973270
$ws.Range("N136").Value = -94806
$ws = $wb.Worksheets.Item("CUL")
# Row 10
$ws.Range("H10").Value = 706.6667
$ws.Range("I10").Value = 65
$ws.Range("J10").Value = 1990
$ws.Range("K10").Value = 195
$ws.Range("L10").Value = 5970
$ws.Range("M10").Value = -56
$ws.Range("N10").Value = -6248
# Row 74
$ws.Range("H74").Value = 8400
$ws.Range("J74").Value = 9480
$ws.Range("L74").Value = 28440
$ws.Range("N74").Value = -30562
# Row 77
$ws.Range("H77").Value = 8400
$ws.Range("J77").Value = 9480
$ws.Range("L77").Value = 85320
$ws.Range("N77").Value = -95928
# Row 104
$ws.Range("H104").Value = 5596.8184
$ws.Range("I104").Value = 1747.5
$ws.Range("J104").Value = 7796.4287
$ws.Range("K104").Value = 5242.5
$ws.Range("L104").Value = 23389.2861
$ws.Range("M104").Value = -2621.5
$ws.Range("N104").Value = -28631.2861
$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 2025.6666
$ws.Range("I126").Value = 1431.1
$ws.Range("J126").Value = 2566.182
$ws.Range("K126").Value = 4293.299999999999
$ws.Range("L126").Value = 7698.545999999999
$ws.Range("M126").Value = -1823.299999999999
$ws.Range("N126").Value = -12638.546
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 1517.3889
$ws.Range("I7").Value = 1322
$ws.Range("J7").Value = 1712.7778
$ws.Range("K7").Value = 1322
$ws.Range("L7").Value = 1712.7778
$ws.Range("M7").Value = -1210
$ws.Range("N7").Value = -1936.7778
# Row 100
$ws.Range("H100").Value = 2158.3333
$ws.Range("I100").Value = 1950
$ws.Range("J100").Value = 2575
$ws.Range("K100").Value = 1950
$ws.Range("L100").Value = 2575
$ws.Range("M100").Value = -1409
$ws.Range("N100").Value = -3657
# Row 126
$ws.Range("H126").Value = 1517.3889
$ws.Range("I126").Value = 1322
$ws.Range("J126").Value = 1712.7778
$ws.Range("K126").Value = 3966
$ws.Range("L126").Value = 5138.3334
$ws.Range("M126").Value = -1496
$ws.Range("N126").Value = -10078.3334
$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 368.1
$ws.Range("I100").Value = 298.5
$ws.Range("J100").Value = 530.5
$ws.Range("K100").Value = 597
$ws.Range("L100").Value = 1061
$ws.Range("M100").Value = -56
$ws.Range("N100").Value = -2143
# Row 122
$ws.Range("H122").Value = 5988.1055
$ws.Range("I122").Value = 7805.2856
$ws.Range("K122").Value = 23415.8568
$ws.Range("M122").Value = -20965.8568
# Row 126
$ws.Range("H126").Value = 173077550
$ws.Range("I126").Value = 138889420
$ws.Range("J126").Value = 250000820
$ws.Range("K126").Value = 416668260
$ws.Range("L126").Value = 750002460
$ws.Range("M126").Value = -416665790
$ws.Range("N126").Value = -750007400
# Row 136
$ws.Range("H136").Value = 22529730
$ws.Range("I136").Value = 12279507
$ws.Range("J136").Value = 55558224
$ws.Range("K136").Value = 36838521
$ws.Range("L136").Value = 166674672
$ws.Range("M136").Value = -36835971
$ws.Range("N136").Value = -166679772
